$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "quantity"

$ws.Range("B6").Select()
